$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Clear()
$ws.Range("F2").Value = "peepeepoopoo123"
